# Applies the Fri Oct  6 17:29:49 UTC 2023 "cryptos list" refresh:
# updated prices / 1h volume deltas, plus a couple of row swaps/replacements
# (PaxDollar <-> WEMIXToken rows 42/43, BabyDogeCoin -> Cronos row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.888.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.55%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.641.19'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.54%  '

# Row 4: TetherUSD
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.21%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.99'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.95%  '

# Row 6: XRP
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.05%  '

# Row 7: USDC
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.26%  '

# Row 8: Solana
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.54'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.68%  '

# Row 9: Cardano
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.26%  '

# Row 10: Dogecoin
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.33%  '

# Row 11: TRON
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.73%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.873.62'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.60%  '

# Row 13: WrappedEther
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.632.75'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.96%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.08'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.29%  '

# Row 15: Polygon
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.51%  '

# Row 16: Litecoin
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.67'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.19%  '

# Row 17: WrappedBTC
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.872.52'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.46%  '

# Row 18: BitcoinCash
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.87'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.83%  '

# Row 19: Chainlink
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.01%  '

# Row 20: ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.63%  '

# Row 21: Dai
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.14%  '

# Row 22: Avalanche
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +8.37%  '

# Row 23: Uniswap
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.35%  '

# Row 24: Toncoin
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.28%  '

# Row 25: Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.50'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.58%  '

# Row 26: Cosmos
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.93'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.48%  '

# Row 27: EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.71'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.99%  '

# Row 28: Stellar
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.06%  '

# Row 29: BinanceUSD
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.22%  '

# Row 30: PancakeSwap
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.99%  '

# Row 31: Hedera
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.29%  '

# Row 32: Filecoin
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.40%  '

# Row 33: Maker
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.458.26'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.95%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.35%  '

# Row 35: LidoDAOToken
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.07%  '

# Row 36: HuobiToken
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.71%  '

# Row 37: ARBITRUM
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.889'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.42%  '

# Row 38: ImmutableX
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.25%  '

# Row 39: VeChain
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.98%  '

# Row 40: TrustWalletToken
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.916'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.60%  '

# Row 41: Aave
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.36'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.25%  '

# Row 42: PaxDollar
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.02'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.82%  '

# Row 43: WEMIXToken
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.20%  '

# Row 44: mCoin
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.16%  '

# Row 45: MXToken
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.64%  '

# Row 46: FraxShare
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.06%  '

# Row 47: RenderToken
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.78'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.20%  '

# Row 48: RocketPoolETH
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.782.60'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.36%  '

# Row 49: Quant
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.51'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.79%  '

# Row 50: Algorand
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.11%  '

# Row 51: BabyDogeCoin
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0507'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.25%  '
